$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G (header "K") for rows 2-14
$values = @{
    2 = 1
    3 = 4
    4 = 2
    5 = 1
    6 = 1
    7 = 1
    8 = 3
    9 = 2
    10 = 4
    11 = 3
    12 = 0
    13 = 0
    14 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
